$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Cells.Item(11, 2).Value = "b2b6ea8d6f2fd751653d2922bf86b7f7"  # 05-050301A
$ws.Cells.Item(24, 2).Value = "2a8524da19a261ecdef6891100f68859"  # 05-050316TC
$ws.Cells.Item(34, 2).Value = "1a2aad99247432a7c8ad2c855eaeec1e"  # 05-050316TP
$ws.Cells.Item(89, 2).Value = "f7945b435d376f43969ae850a7cc68cb"  # 05-050104A
$ws.Cells.Item(99, 2).Value = "45345d734b099da46e786c83e8f28c96"  # 05-050101A
$ws.Cells.Item(110, 2).Value = "74c498ae62afc36eaf69fb2be262b624"  # 05-050102A
$ws.Cells.Item(121, 2).Value = "16e942b2f0271e54d831782a253ff8bb"  # 05-050301TP
$ws.Cells.Item(154, 2).Value = "6b15316edc1cc092b4abac42be90bd28"  # 05-050007TC
$ws.Cells.Item(160, 2).Value = "a971ea9eb8c3823f3586968e3793190b"  # 05-050007TP
$ws.Cells.Item(162, 2).Value = "b2958ca0a2f48c38ed413b0942283382"  # 05-050308A
$ws.Cells.Item(175, 2).Value = "d7a63990157d9dcc566c9c52d107a4bf"  # 05-050303TP
$ws.Cells.Item(180, 2).Value = "9ff250cc2296e8b04e2e9c55eb7b492a"  # 05-050303TC
$ws.Cells.Item(191, 2).Value = "7a37b708a3b20888dceb268fa7491c0d"  # 05-050314TP
$ws.Cells.Item(213, 2).Value = "289d9c7f686850f0271f00b042591a5a"  # 05-050303A
$ws.Cells.Item(281, 2).Value = "beba7bce29c4068483cd10898052ff4a"  # 05-050101TP
$ws.Cells.Item(338, 2).Value = "7d3192fea74a6be1ead9e53c83c35f0f"  # 05-050005TP
$ws.Cells.Item(488, 2).Value = "238ad20f0552c8e5601f2bcb733f8eba"  # 05-050314A
$ws.Cells.Item(511, 2).Value = "3bb24bf20af84bd73d4fd48e30da03f3"  # 05-050208TP
$ws.Cells.Item(516, 2).Value = "3573f972709eca56275fd504bb286c75"  # 05-050306TP
$ws.Cells.Item(524, 2).Value = "e3d6f2571a6e47a237de56acc60583d0"  # 05-050317TC
$ws.Cells.Item(535, 2).Value = "c2ff6a83c1beba8689e2d6eaa3eb06e1"  # 05-050317TP
$ws.Cells.Item(545, 2).Value = "caed40e30b8d326c9ee29159f49801d9"  # 05-050006A
$ws.Cells.Item(559, 2).Value = "94c8a699ba72fa2ba49483e62eaeeb5b"  # 05-050201A
$ws.Cells.Item(565, 2).Value = "6dae6fa19d878e3e786208dc34f13627"  # 05-050007A
$ws.Cells.Item(596, 2).Value = "0500c3294f2fe90971052abfee60871b"  # 05-050005A
$ws.Cells.Item(677, 2).Value = "2ede366eee4394e48ea0925f9464345c"  # 05-050208A
$ws.Cells.Item(678, 2).Value = "7f37c26eae181fa0ad2e97b5864751b2"  # 05-050317A
$ws.Cells.Item(716, 2).Value = "14cb8d34718c47516b19ad2970bcf17c"  # 05-050315A
$ws.Cells.Item(741, 2).Value = "93049bfcc2ff1ccbc37fcd3a7fe75f92"  # 05-050316A
$ws.Cells.Item(754, 2).Value = "73dcb4033cf74069e3da205ee99500a5"  # 05-050315TP
$ws.Cells.Item(780, 2).Value = "0a647b4a3f32e50bca26867df944df5e"  # 05-050102TP
$ws.Cells.Item(823, 2).Value = "d05f60cb7fe7ed68b218c83ac767a514"  # 05-050006TC
$ws.Cells.Item(827, 2).Value = "828dfcdbe017b46b27ba6a91372baea2"  # 05-050006TP
$ws.Cells.Item(828, 2).Value = "369163dccc3c430a954a07963037cfd1"  # 05-050104TC
$ws.Cells.Item(837, 2).Value = "55ee70e9919cf8142a528225a340560d"  # 05-050104TM
$ws.Cells.Item(839, 2).Value = "e8dfad8ff97156163b1440cb8b6475c6"  # 05-050104TP
$ws.Cells.Item(886, 2).Value = "d878f735a89572d2273c1e98708e28dd"  # 03-030032A
